$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vacation Summary")
$gs = $wb.Worksheets.Item("General Summary")

# --- Vacation Summary: insert a new "DI Apprentice" row above the existing
#     "Sup I" / "Sup II" rows (rows 3 & 4), pushing the existing data down
#     by one logical row without changing the sheet's physical row count
#     (row 5 was already a formatted-but-empty row, so it absorbs the shift).

# 1) Capture the current (pre-edit) row 3 and row 4 data before it gets
#    overwritten.
$oldA3 = $ws.Range("A3").Value2
$oldB3 = $ws.Range("B3").Value2
$oldC3 = $ws.Range("C3").Value2
$oldF3 = $ws.Range("F3").Value2

$oldA4 = $ws.Range("A4").Value2
$oldB4 = $ws.Range("B4").Value2
$oldC4 = $ws.Range("C4").Value2
$oldF4 = $ws.Range("F4").Value2

# 2) Row 5 (A:H) has no formatting yet (it was a blank placeholder row) -
#    copy the formatting from row 4 so the new content matches the rest
#    of the table.
$ws.Range("A4:H4").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 3) Push the old row 4 ("Sup II") data down into row 5.
$ws.Range("A5").Value = $oldA4
$ws.Range("B5").Value = $oldB4
$ws.Range("C5").Value = $oldC4
$ws.Range("F5").Value = $oldF4
$ws.Range("D5").Formula = "=C5/7.5"
$ws.Range("E5").Formula = "=C5/1957.5"
$ws.Range("G5").Formula = "='General Summary'!`$C`$2-F5"
$ws.Range("H5").Formula = "=F5/'General Summary'!`$C`$2"

# 4) Push the old row 3 ("Sup I") data down into row 4.
$ws.Range("A4").Value = $oldA3
$ws.Range("B4").Value = $oldB3
$ws.Range("C4").Value = $oldC3
$ws.Range("F4").Value = $oldF3
$ws.Range("D4").Formula = "=C4/7.5"
$ws.Range("E4").Formula = "=C4/1957.5"
$ws.Range("G4").Formula = "='General Summary'!`$C`$2-F4"
$ws.Range("H4").Formula = "=F4/'General Summary'!`$C`$2"

# 5) Write the new "DI Apprentice" row into row 3 (style already present).
$ws.Range("A3").Value = "DI Apprentice"
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = $oldC3
$ws.Range("F3").Value = $oldF3
$ws.Range("D3").Formula = "=C3/7.5"
$ws.Range("E3").Formula = "=C3/1957.5"
$ws.Range("G3").Formula = "='General Summary'!`$C`$2-F3"
$ws.Range("H3").Formula = "=F3/'General Summary'!`$C`$2"

# --- Sheet-view / active-tab bookkeeping: the "Vacation Summary" tab is now
#     the active one (was "General Summary" before), with a new selection.
$gs.Activate()
$gs.Range("C4").Select()
$ws.Activate()
$ws.Range("K28").Select()
